$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (ID "H 72") entirely; this shifts every
# subsequent row up by one and correctly updates the used range.
$ws.Rows.Item(2).Delete()

# The B/C/E ("A"/"B"/"D" header) columns hold a randomized missing-value
# mask for this seed; re-apply the new mask/values for each data row.
$ws.Range("B2").Value = -20
$ws.Range("C2").Value = 10
$ws.Range("E2").Value = -6.9
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 11.7
$ws.Range("E3").ClearContents()
$ws.Range("B4").Value = -19.7
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("B5").Value = -19.7
$ws.Range("C5").Value = 11.2
$ws.Range("E5").Value = -5.7
$ws.Range("B6").Value = -18.7
$ws.Range("C6").Value = 11
$ws.Range("E6").ClearContents()
$ws.Range("B7").Value = -19.5
$ws.Range("C7").Value = 12.3
$ws.Range("E7").Value = -5
$ws.Range("B8").Value = -19.8
$ws.Range("C8").ClearContents()
$ws.Range("E8").Value = -5.7
$ws.Range("B9").Value = -19.9
$ws.Range("C9").Value = 12
$ws.Range("E9").ClearContents()
$ws.Range("B10").Value = -19.5
$ws.Range("C10").Value = 15
$ws.Range("E10").ClearContents()
$ws.Range("B11").Value = -19.9
$ws.Range("C11").Value = 15.5
$ws.Range("E11").Value = -6.6
$ws.Range("B12").Value = -20.6
$ws.Range("C12").Value = 10.5
$ws.Range("E12").Value = -6.8
$ws.Range("B13").Value = -19.8
$ws.Range("C13").Value = 11.5
$ws.Range("E13").Value = -6.1
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = 11.4
$ws.Range("E14").Value = -7.9
$ws.Range("B15").Value = -18.9
$ws.Range("C15").Value = 12
$ws.Range("E15").Value = -12
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 12.5
$ws.Range("E16").Value = -5.3
$ws.Range("B17").Value = -19.9
$ws.Range("C17").Value = 12.5
$ws.Range("E17").Value = -5.3
$ws.Range("B18").Value = -19.6
$ws.Range("C18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("B19").Value = -19.1
$ws.Range("C19").Value = 12.5
$ws.Range("E19").Value = -8.4
$ws.Range("B20").Value = -19.5
$ws.Range("C20").Value = 13.5
$ws.Range("E20").Value = -6.9
$ws.Range("B21").ClearContents()
$ws.Range("C21").Value = 11.2
$ws.Range("E21").ClearContents()
$ws.Range("B22").Value = -19.6
$ws.Range("C22").ClearContents()
$ws.Range("E22").Value = -8.5
$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = 13.2
$ws.Range("E23").Value = -6.5
$ws.Range("B24").Value = -19
$ws.Range("C24").Value = 12.5
$ws.Range("E24").Value = -7.2
$ws.Range("B25").ClearContents()
$ws.Range("C25").Value = 12.7
$ws.Range("E25").Value = -8.699999999999999
$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = 12
$ws.Range("E26").Value = -8.5
$ws.Range("B27").Value = -19.3
$ws.Range("C27").Value = 13.5
$ws.Range("E27").ClearContents()
$ws.Range("B28").Value = -19.5
$ws.Range("C28").Value = 12.2
$ws.Range("E28").Value = -7
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 12.7
$ws.Range("E29").Value = -8.1
$ws.Range("B30").Value = -19.5
$ws.Range("C30").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("B31").Value = -19.5
$ws.Range("C31").Value = 10.7
$ws.Range("E31").ClearContents()
$ws.Range("B32").Value = -19.4
$ws.Range("C32").ClearContents()
$ws.Range("E32").Value = -7.8
$ws.Range("B33").Value = -19.4
$ws.Range("C33").Value = 12
$ws.Range("E33").Value = -9.5
$ws.Range("B34").Value = -19.6
$ws.Range("C34").Value = 10.3
$ws.Range("E34").ClearContents()
$ws.Range("B35").Value = -19.2
$ws.Range("C35").Value = 11.3
$ws.Range("E35").ClearContents()
$ws.Range("B36").Value = -19.1
$ws.Range("C36").Value = 14.3
$ws.Range("E36").ClearContents()
$ws.Range("B37").Value = -19.8
$ws.Range("C37").Value = 12.1
$ws.Range("E37").Value = -7.1
$ws.Range("B38").Value = -19.2
$ws.Range("C38").Value = 11.5
$ws.Range("E38").Value = -9.199999999999999
$ws.Range("B39").Value = -19.8
$ws.Range("C39").ClearContents()
$ws.Range("E39").ClearContents()
$ws.Range("B40").ClearContents()
$ws.Range("C40").Value = 11.9
$ws.Range("E40").Value = -7.9
$ws.Range("B41").Value = -18.1
$ws.Range("C41").Value = 13.9
$ws.Range("E41").Value = -10.2
$ws.Range("B42").Value = -19
$ws.Range("C42").ClearContents()
$ws.Range("E42").Value = -6.8
$ws.Range("B43").Value = -18.9
$ws.Range("C43").Value = 13
$ws.Range("E43").Value = -8.800000000000001
$ws.Range("B44").Value = -18.8
$ws.Range("C44").Value = 12.6
$ws.Range("E44").Value = -8.9
$ws.Range("B45").Value = -19.7
$ws.Range("C45").Value = 11.7
$ws.Range("E45").Value = -7.4
$ws.Range("B46").Value = -19.5
$ws.Range("C46").Value = 11.7
$ws.Range("E46").Value = -7.3
$ws.Range("B47").Value = -19.8
$ws.Range("C47").Value = 12.9
$ws.Range("E47").Value = -6.6
$ws.Range("B48").Value = -19.3
$ws.Range("C48").Value = 12.7
$ws.Range("E48").Value = -8.300000000000001
$ws.Range("B49").Value = -19.9
$ws.Range("C49").Value = 11.5
$ws.Range("E49").Value = -7.2
$ws.Range("B50").Value = -19.7
$ws.Range("C50").Value = 10.7
$ws.Range("E50").Value = -8.800000000000001
$ws.Range("B51").Value = -20.5
$ws.Range("C51").Value = 11.6
$ws.Range("E51").Value = -7.7
$ws.Range("B52").Value = -20.2
$ws.Range("C52").Value = 10.8
$ws.Range("E52").ClearContents()
$ws.Range("B53").ClearContents()
$ws.Range("C53").Value = 10.5
$ws.Range("E53").ClearContents()
$ws.Range("B54").Value = -17.2
$ws.Range("C54").Value = 14.3
$ws.Range("E54").Value = -6.3
$ws.Range("B55").Value = -20.4
$ws.Range("C55").Value = 10
$ws.Range("E55").Value = -10
$ws.Range("B56").Value = -19.2
$ws.Range("C56").Value = 11.9
$ws.Range("E56").Value = -5.7
$ws.Range("B57").ClearContents()
$ws.Range("C57").Value = 11.1
$ws.Range("E57").Value = -5.9
$ws.Range("B58").Value = -19.5
$ws.Range("C58").ClearContents()
$ws.Range("E58").Value = -6.8
$ws.Range("B59").ClearContents()
$ws.Range("C59").Value = 11.4
$ws.Range("E59").Value = -5.7
$ws.Range("B60").Value = -18.8
$ws.Range("C60").Value = 15.3
$ws.Range("E60").Value = -8.1
$ws.Range("B61").Value = -19.9
$ws.Range("C61").Value = 10.5
$ws.Range("E61").Value = -6.4
$ws.Range("B62").Value = -19.5
$ws.Range("C62").Value = 10.4
$ws.Range("E62").Value = -10.7
